$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 0. Update the first title (row 1 text only - style/merge already correct)
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "Location Aspect Type template"

# ---------------------------------------------------------------------------
# 1. Remove the old second table (rows 18-34) - unmerge, then clear.
#    Rows 1-14 (first table) stay untouched and double as style "seed" rows
#    for copy/paste-format operations below (A1/B1 = title style,
#    A2/B2 = sub-header style, B3 = red-highlighted-value style,
#    A3 = plain default style).
# ---------------------------------------------------------------------------
$ws.Range("A18:B18").UnMerge()
$ws.Range("A16:B34").Clear()

# ---------------------------------------------------------------------------
# 2. New section: "Location Aspect Type Attributes" (rows 16-26)
# ---------------------------------------------------------------------------

# Row 16 - section title, merged
$ws.Range("A1:B1").Copy()
$ws.Range("A16:B16").PasteSpecial(-4122)
$ws.Range("A16").Value = "Location Aspect Type Attributes"
$ws.Range("A16:B16").Merge()

# Row 17 - sub-header (Attribute / (multi)Selection Values)
$ws.Range("A2:B2").Copy()
$ws.Range("A17:B17").PasteSpecial(-4122)
$ws.Range("A17").Value = "Attribute"
$ws.Range("B17").Value = "(multi)Selection Values"
$ws.Rows("17:17").RowHeight = 24.75

# Row 18
$ws.Range("A18").Value = "Hazardous classification req."
$ws.Range("B3").Copy()
$ws.Range("B18").PasteSpecial(-4122)
$ws.Range("B18").Value = "Safe by ventilation, Safe by location, Zone 2, Zone 1, Zone 0"

# Row 19 (plain style on both columns)
$ws.Range("A19").Value = "Fire protection requirements"
$ws.Range("B19").Value = "Cladding, Portable, Sprinkler, Deluge, Foam"

# Row 20 (plain style on both columns)
$ws.Range("A20").Value = "Fire detection requirements"
$ws.Range("B20").Value = "Call points, Gas, Flame, Heat, Smoke"

# Row 21
$ws.Range("A21").Value = "Environment classification req."
$ws.Range("B3").Copy()
$ws.Range("B21").PasteSpecial(-4122)
$ws.Range("B21").Value = "Heated, Cooled & Ventilated, Cooled & Ventilated, Weather Protected, Weather Exposed"

# Row 22
$ws.Range("A22").Value = "Noise limitation requirements"
$ws.Range("B3").Copy()
$ws.Range("B22").PasteSpecial(-4122)
$ws.Range("B22").Value = "55dB, 70dB, 85dB"

# Row 23 (plain style on both columns)
$ws.Range("A23").Value = "Wireless coverage requrements"
$ws.Range("B23").Value = "Wifi; UHF; PRS"

# Row 24 intentionally blank (separator)

# Row 25
$ws.Range("A25").Value = "Area activity classification"
$ws.Range("B3").Copy()
$ws.Range("B25").PasteSpecial(-4122)
$ws.Range("B25").Value = "Process, HVAC, Electrical, Control, Laboratory, Workshop, Access, Trucking"

# Row 26 (plain style on both columns)
$ws.Range("A26").Value = "Area partition classifications"
$ws.Range("B26").Value = "Table"

# ---------------------------------------------------------------------------
# 3. New section: "Location Aspect Instance Attributes" (rows 28-34)
# ---------------------------------------------------------------------------

# Row 28 - section title, merged
$ws.Range("A1:B1").Copy()
$ws.Range("A28:B28").PasteSpecial(-4122)
$ws.Range("A28").Value = "Location Aspect Instance Attributes"
$ws.Range("A28:B28").Merge()

# Row 29 - sub-header (Attribute / (multi)Selection Values)
$ws.Range("A2:B2").Copy()
$ws.Range("A29:B29").PasteSpecial(-4122)
$ws.Range("A29").Value = "Attribute"
$ws.Range("B29").Value = "(multi)Selection Values"

# Row 30
$ws.Range("A30").Value = "Relative North"
$ws.Range("B3").Copy()
$ws.Range("B30").PasteSpecial(-4122)
$ws.Range("B30").Value = "mm"

# Row 31
$ws.Range("A31").Value = "Relative East"
$ws.Range("B3").Copy()
$ws.Range("B31").PasteSpecial(-4122)
$ws.Range("B31").Value = "mm"

# Row 32
$ws.Range("A32").Value = "Relative Elevation"
$ws.Range("B3").Copy()
$ws.Range("B32").PasteSpecial(-4122)
$ws.Range("B32").Value = "mm"

# Row 33
$ws.Range("A33").Value = "Size box WHD"
$ws.Range("B3").Copy()
$ws.Range("B33").PasteSpecial(-4122)
$ws.Range("B33").Value = "m x m x m OR mm x mm x mm"

# Row 34
$ws.Range("A34").Value = "Depth BSL (subsea only)"
$ws.Range("B3").Copy()
$ws.Range("B34").PasteSpecial(-4122)
$ws.Range("B34").Value = "m"

$ws.Application.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 4. New explanatory note (row 37), wrapped text, merged
# ---------------------------------------------------------------------------
$note = "A Location Type is constructed by selecting main type/sub-type, then charactising the type further by setting type attributes.`nE.g. Attribute values: Area/Sub Area, Safe by location, Deluge, Call points+Flame, Weather Protected, 85dB, UHF+PRS, Process`nwold define a location type representing where the pumps are located. When the Type is instansiated in the model, the Instance Attributes are given values - such as coordinates and size box."

$ws.Range("A37").Value = $note
$ws.Range("A37").HorizontalAlignment = -4131
$ws.Range("A37").VerticalAlignment = -4160
$ws.Range("A37").WrapText = $true

$ws.Range("B37").HorizontalAlignment = -4131
$ws.Range("B37").VerticalAlignment = -4160

$ws.Range("A37:B37").Merge()
$ws.Rows("37:37").RowHeight = 88.5

# ---------------------------------------------------------------------------
# 5. Selection (matches final saved cursor position in the diff)
# ---------------------------------------------------------------------------
$ws.Range("B42").Select()
